$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new dividend row (row 84) below the existing data (row 83),
# carrying over row 83's formatting (the date-number style on column A)
# the same way Excel's "fill the row below" / copy-paste workflow would.
$ws.Range("A83:B83").Copy()
$ws.Range("A84").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A84").Value = 45422   # 2024-05-10
$ws.Range("B84").Value = 0.25
